$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper constants (PasteSpecial paste type for "Formats only")
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# New data rows (20-33) appended below the existing table (A1:E18).
# Column A = KPI Set, B = KPI Name, C = Atomic Name Old, D = Atomic Name New,
# E = generated SQL (CONCATENATE formula identical in shape to existing rows).
# ---------------------------------------------------------------------------

function Set-UpdateFormula($row) {
    $formula = '=CONCATENATE("UPDATE `static`.atomic_kpi a JOIN `static`.kpi k ON k.pk=a.kpi_fk JOIN `static`.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name=''",D' + $row + ',"'', a.description=''",D' + $row + ',"'', a.display_text=''",D' + $row + ',"''  WHERE s.name=''",A' + $row + ',"'' AND k.display_text=''",B' + $row + ',"'' AND a.name=''",C' + $row + ',"'';")'
    $ws.Range("E$row").Formula = $formula
}

# Row 20
$ws.Range("A20").Value = "PoS 2019 - IC Petroleum - CAP"
$ws.Range("B20").Value = "Juice Availability"
$ws.Range("C20").Value = "NEW SKU 7"
$ws.Range("D20").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A2").Copy()
$ws.Range("A20:C20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# build the brand new highlighted style (green fill / teal pattern colour)
# on D20 - base it on the plain-bordered style (B11) then recolour the
# interior so it ends up matching the "new SKU" highlight used below.
$ws.Range("B11").Copy()
$ws.Range("D20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("D20").Interior.Color = 5287936
$ws.Range("D20").Interior.PatternColor = 8421376

Set-UpdateFormula 20

# Row 21
$ws.Range("A21").Value = "PoS 2019 - IC Petroleum - CAP"
$ws.Range("B21").Value = "Juice Availability"
$ws.Range("C21").Value = "NEW SKU 8"
$ws.Range("D21").Value = "Pulpy - Watermelon-Strawberry - 0.45L"

$ws.Range("A2").Copy()
$ws.Range("A21:C21").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 21

# Row 22
$ws.Range("A22").Value = "PoS 2019 - IC Petroleum - CAP"
$ws.Range("B22").Value = "Juice Availability"
$ws.Range("C22").Value = "NEW SKU 7"
$ws.Range("D22").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A2").Copy()
$ws.Range("A22:C22").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("D22").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 22

# Row 23
$ws.Range("A23").Value = "PoS 2019 - IC Petroleum - CAP"
$ws.Range("B23").Value = "Juice Availability"
$ws.Range("C23").Value = "NEW SKU 8"
$ws.Range("D23").Value = "Pulpy - Watermelon-Strawberry - 0.45L"

$ws.Range("A2").Copy()
$ws.Range("A23:C23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("D23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 23

# Row 24
$ws.Range("A24").Value = "PoS 2019 - MT Conv Big - CAP"
$ws.Range("B24").Value = "Juice (JNSD) Availability"
$ws.Range("C24").Value = "NEW SKU 6"
$ws.Range("D24").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A24").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("B24").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C24:D24").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 24

# Row 25
$ws.Range("A25").Value = "PoS 2019 - MT Conv Big - REG"
$ws.Range("B25").Value = "Juice (JNSD) Availability"
$ws.Range("C25").Value = "NEW SKU 6"
$ws.Range("D25").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A25").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("B25").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C25:D25").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 25

# Row 26
$ws.Range("A26").Value = "PoS 2019 - MT Hypermarket - CAP"
$ws.Range("B26").Value = "Juice (JNSD) Availability"
$ws.Range("C26").Value = "NEW SKU 6"
$ws.Range("D26").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A26").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B26").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C26:D26").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 26

# Row 27
$ws.Range("A27").Value = "PoS 2019 - MT Hypermarket - CAP"
$ws.Range("B27").Value = "Juice (JNSD) Availability"
$ws.Range("C27").Value = "NEW SKU 8"
$ws.Range("D27").Value = "Pulpy - Watermelon-Strawberry - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C27:D27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 27

# Row 28
$ws.Range("A28").Value = "PoS 2019 - MT Hypermarket - REG"
$ws.Range("B28").Value = "Juice (JNSD) Availability"
$ws.Range("C28").Value = "NEW SKU 6"
$ws.Range("D28").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C28:D28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 28

# Row 29
$ws.Range("A29").Value = "PoS 2019 - MT Hypermarket - REG"
$ws.Range("B29").Value = "Juice (JNSD) Availability"
$ws.Range("C29").Value = "NEW SKU 8"
$ws.Range("D29").Value = "Pulpy - Watermelon-Strawberry - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A29").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B29").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C29:D29").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 29

# Row 30
$ws.Range("A30").Value = "PoS 2019 - MT Supermarket - CAP"
$ws.Range("B30").Value = "Juice (JNSD) Availability"
$ws.Range("C30").Value = "NEW SKU 6"
$ws.Range("D30").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A30").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B30").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C30:D30").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 30

# Row 31
$ws.Range("A31").Value = "PoS 2019 - MT Supermarket - CAP"
$ws.Range("B31").Value = "Juice (JNSD) Availability"
$ws.Range("C31").Value = "NEW SKU 8"
$ws.Range("D31").Value = "Pulpy - Watermelon-Strawberry - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A31").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B31").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C31:D31").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 31

# Row 32
$ws.Range("A32").Value = "PoS 2019 - MT Supermarket - REG"
$ws.Range("B32").Value = "Juice (JNSD) Availability"
$ws.Range("C32").Value = "NEW SKU 6"
$ws.Range("D32").Value = "Pulpy - Guava-Passion Fruit - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A32").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B32").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C32:D32").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 32

# Row 33
$ws.Range("A33").Value = "PoS 2019 - MT Supermarket - REG"
$ws.Range("B33").Value = "Juice (JNSD) Availability"
$ws.Range("C33").Value = "NEW SKU 8"
$ws.Range("D33").Value = "Pulpy - Watermelon-Strawberry - 0.45L"

$ws.Range("A18").Copy()
$ws.Range("A33").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("B11").Copy()
$ws.Range("B33").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("C33:D33").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

Set-UpdateFormula 33

# ---------------------------------------------------------------------------
# Row heights for the newly added rows match the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A20:A33").EntireRow.RowHeight = 13.8

# ---------------------------------------------------------------------------
# Update the sheet selection so it lands on the newly added block, mirroring
# what Excel stores after a user finishes editing these new rows.
# ---------------------------------------------------------------------------
$ws.Range("E20:E33").Select()

Write-Host "Applied KPI atomic name rows 20-33"
